$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh -- cell-by-cell update matching the source diff.
# Number-looking Price (column D) values are written with a leading apostrophe
# (forces text entry, same as typing it in Excel) and the cell Style is restored
# immediately after so no stray NumberFormat/quotePrefix style sticks around.

$ws.Cells.Item(2, 4).Value = '69.965.94'
$ws.Cells.Item(2, 5).Value = '  +0.40%  '
$ws.Cells.Item(3, 4).Value = '3.781.81'
$ws.Cells.Item(3, 5).Value = '  +3.63%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).Value = "'615.64"
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  +4.34%  '
$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).Value = "'178.11"
$ws.Cells.Item(6, 4).Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  -4.25%  '
$ws.Cells.Item(7, 4).Value = '3.778.61'
$ws.Cells.Item(7, 5).Value = '  +3.62%  '
$ws.Cells.Item(8, 5).Value = '  +0.10%  '
$origStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).Value = "'0.535"
$ws.Cells.Item(9, 4).Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  +0.28%  '
$ws.Cells.Item(10, 5).Value = '  +3.87%  '
$ws.Cells.Item(11, 5).Value = '  -3.16%  '
$origStyle = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).Value = "'0.493"
$ws.Cells.Item(12, 4).Style = $origStyle
$ws.Cells.Item(12, 5).Value = '  -0.61%  '
$origStyle = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).Value = "'41.05"
$ws.Cells.Item(13, 4).Style = $origStyle
$ws.Cells.Item(13, 5).Value = '  +4.29%  '
$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).Value = "'0.0000255"
$ws.Cells.Item(14, 4).Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  +0.98%  '
$ws.Cells.Item(15, 4).Value = '4.412.90'
$ws.Cells.Item(15, 5).Value = '  +3.57%  '
$ws.Cells.Item(16, 4).Value = '3.784.93'
$ws.Cells.Item(16, 5).Value = '  +3.66%  '
$ws.Cells.Item(17, 4).Value = '70.007.01'
$ws.Cells.Item(17, 5).Value = '  +0.16%  '
$ws.Cells.Item(18, 5).Value = '  -0.07%  '
$origStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).Value = "'7.60"
$ws.Cells.Item(19, 4).Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  +0.90%  '
$origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).Value = "'515.32"
$ws.Cells.Item(20, 4).Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  +1.14%  '
$origStyle = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).Value = "'16.63"
$ws.Cells.Item(21, 4).Style = $origStyle
$ws.Cells.Item(21, 5).Value = '  -3.75%  '
$origStyle = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).Value = "'9.59"
$ws.Cells.Item(22, 4).Style = $origStyle
$ws.Cells.Item(22, 5).Value = '  +2.53%  '
$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).Value = "'0.728"
$ws.Cells.Item(23, 4).Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  -2.66%  '
$ws.Cells.Item(24, 5).Value = '  +5.33%  '
$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).Value = "'87.98"
$ws.Cells.Item(25, 4).Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  +0.04%  '
$origStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).Value = "'13.34"
$ws.Cells.Item(26, 4).Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  -1.65%  '
$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).Value = "'11.09"
$ws.Cells.Item(27, 4).Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  +2.00%  '
$ws.Cells.Item(28, 5).Value = '  +23.92%  '
$ws.Cells.Item(29, 5).Value = '  +0.00%  '
$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).Value = "'2.49"
$ws.Cells.Item(30, 4).Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  -2.28%  '
$ws.Cells.Item(31, 2).Value = 'NEARProtocol'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).Value = "'7.83"
$ws.Cells.Item(31, 4).Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  -4.48%  '
$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).Value = "'2.83"
$ws.Cells.Item(32, 4).Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  +3.18%  '
$origStyle = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).Value = "'31.72"
$ws.Cells.Item(33, 4).Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  -2.35%  '
$origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).Value = "'0.115"
$ws.Cells.Item(34, 4).Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  -1.67%  '
$origStyle = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).Value = "'1.00"
$ws.Cells.Item(35, 4).Style = $origStyle
$ws.Cells.Item(35, 5).Value = '  +0.13%  '
$origStyle = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).Value = "'6.23"
$ws.Cells.Item(36, 4).Style = $origStyle
$ws.Cells.Item(36, 5).Value = '  +1.07%  '
$ws.Cells.Item(37, 5).Value = '  +2.56%  '
$ws.Cells.Item(38, 5).Value = '  +1.70%  '
$origStyle = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).Value = "'2.17"
$ws.Cells.Item(39, 4).Style = $origStyle
$ws.Cells.Item(39, 5).Value = '  +2.74%  '
$origStyle = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).Value = "'0.132"
$ws.Cells.Item(40, 4).Style = $origStyle
$ws.Cells.Item(40, 5).Value = '  +3.42%  '
$origStyle = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).Value = "'51.37"
$ws.Cells.Item(41, 4).Style = $origStyle
$ws.Cells.Item(41, 5).Value = '  +1.25%  '
$origStyle = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).Value = "'44.45"
$ws.Cells.Item(42, 4).Style = $origStyle
$ws.Cells.Item(42, 5).Value = '  -4.89%  '
$origStyle = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).Value = "'8.77"
$ws.Cells.Item(43, 4).Style = $origStyle
$ws.Cells.Item(43, 5).Value = '  -0.70%  '
$origStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).Value = "'424.04"
$ws.Cells.Item(44, 4).Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  +4.53%  '
$ws.Cells.Item(45, 4).Value = '3.064.51'
$ws.Cells.Item(45, 5).Value = '  -3.24%  '
$origStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).Value = "'2.74"
$ws.Cells.Item(46, 4).Style = $origStyle
$ws.Cells.Item(46, 5).Value = '  -1.49%  '
$origStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).Value = "'0.0364"
$ws.Cells.Item(47, 4).Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  -0.57%  '
$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).Value = "'27.72"
$ws.Cells.Item(48, 4).Style = $origStyle
$ws.Cells.Item(48, 5).Value = '  -0.39%  '
$ws.Cells.Item(49, 2).Value = 'USDe'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$origStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).Value = "'1.00"
$ws.Cells.Item(49, 4).Style = $origStyle
$ws.Cells.Item(49, 5).Value = '  -0.02%  '
$ws.Cells.Item(50, 2).Value = 'ThetaToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$origStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).Value = "'2.49"
$ws.Cells.Item(50, 4).Style = $origStyle
$ws.Cells.Item(50, 5).Value = '  +1.73%  '
$origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).Value = "'135.23"
$ws.Cells.Item(51, 4).Style = $origStyle
$ws.Cells.Item(51, 5).Value = '  -1.05%  '
